# The deck ships two theme parts:
#   ppt/theme/theme1.xml  -> bound to the slide master (the "Integral" / Red
#                             Violet palette) - this is the theme every slide
#                             actually renders with.
#   ppt/theme/theme2.xml  -> bound to the notes master (plain "Office Theme"
#                             palette).
# The authored change swaps which palette lives in which part: the slide
# master's theme becomes the stock Office palette, while the notes master's
# theme becomes the former Integral/Red-Violet palette. Font scheme and
# format scheme (fills/lines/effects) are identical between the two themes
# already, so only the 12 theme colors differ.
#
# Apply the new (Office) palette to the deck's live theme via the
# ThemeColorScheme object, which is keyed in the fixed order:
#   1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3, 8 accent4,
#   9 accent5, 10 accent6, 11 hlink, 12 folHlink

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

$tcs.Colors(1).RGB  = 0         # dk1      000000
$tcs.Colors(2).RGB  = 16777215  # lt1      FFFFFF
$tcs.Colors(3).RGB  = 6968388   # dk2      44546A
$tcs.Colors(4).RGB  = 15132391  # lt2      E7E6E6
$tcs.Colors(5).RGB  = 13998939  # accent1  5B9BD5
$tcs.Colors(6).RGB  = 3243501   # accent2  ED7D31
$tcs.Colors(7).RGB  = 10855845  # accent3  A5A5A5
$tcs.Colors(8).RGB  = 49407     # accent4  FFC000
$tcs.Colors(9).RGB  = 12874308  # accent5  4472C4
$tcs.Colors(10).RGB = 4697456   # accent6  70AD47
$tcs.Colors(11).RGB = 12673797  # hlink    0563C1
$tcs.Colors(12).RGB = 7491477   # folHlink 954F72
